$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new row of data: name, email, repo link
$ws.Range("A2").Value = "احمد اسماعيل محمود حسن"
$ws.Range("B2").Value = "Ahmedaldarawy3@gmail.com"
$ws.Range("C2").Value = "https://github.com/E0xMomen/open_source_project"

# Turn the email into a mailto hyperlink (this also applies the
# built-in "Hyperlink" style - underline + theme color - to B2)
[void]$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:Ahmedaldarawy3@gmail.com")

# Match the author's final selection
[void]$ws.Range("C4").Select()
